$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 467, shifting existing rows 467-505 down to 468-506.
$ws.Rows.Item(467).Insert()

# Populate the newly inserted row 467 with the new data record.
$ws.Range("A467").Value = 3
$ws.Range("B467").Value = "Femacal de La Calera"
$ws.Range("C467").Value = "Coquimbo"
$ws.Range("D467").Value = 44769
$ws.Range("E467").Value = 5
$ws.Range("F467").Value = 100112032
$ws.Range("G467").Value = "Zapallo italiano"
$ws.Range("H467").Value = "Sin especificar"
$ws.Range("I467").Value = "Primera"
$ws.Range("J467").Value = 235
$ws.Range("K467").Value = 13000
$ws.Range("L467").Value = 14000
$ws.Range("M467").Value = 13213
$ws.Range("N467").Value = "$/caja 70 unidades"
$ws.Range("O467").Value = "Región de Arica y Parinacota"
$ws.Range("P467").Value = 189
$ws.Range("Q467").Value = 70
$ws.Range("R467").Value = "Hortaliza"
